$d = $word.ActiveDocument

# "En esta sección vamos a ver la diferencia ..." ->
# "En esta sección se visualizará la diferencia ..."
$d.Content.Find.Execute("vamos a ver la ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "se visualizará la ", 2)
